$d = $word.ActiveDocument

# --- 1. Title paragraph: ${project.project_name} -> ${project.name}, with a
#        _GoBack bookmark left at the edit point (between "${project." and "name}") ---
$p1 = $d.Paragraphs(1)
$old1 = "`$`{project.project_name}"
$new1 = "`$`{project.name}"
$p1.Range.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# Force the run to split right after "${project." so the bookmark lands between
# the two runs exactly like Word leaves the "last edit" marker.
$old1b = "name}"
$new1b = "name}"
$p1.Range.Find.Execute($old1b, $false, $false, $false, $false, $false, $true, 1, $false, $new1b, 2) | Out-Null

$splitPos = $p1.Range.Start + 10
$goBackRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# --- 2. "non_research_goals" placeholder paragraph: merge "project." into the
#        following run, and mark the paragraph indent with CharacterUnitLeftIndent=0 ---
$p5 = $d.Paragraphs(5)
$old2 = "project.non_research_goals"
$new2 = "project.non_research_goals"
$p5.Range.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null
$p5.Range.ParagraphFormat.CharacterUnitLeftIndent = 0

# --- 3. Remove the stray _GoBack bookmark that previously sat in the trailing
#        empty paragraph (a new one was already placed above). ---
$d.Bookmarks("_GoBack").Delete() | Out-Null
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# --- 4. Add header/footer parts (even / default / first) to the only section ---
$sec = $d.Sections(1)
$sec.PageSetup.DifferentFirstPageHeaderFooter = $true
$sec.Headers(1).Range.Text = ""
$sec.Headers(2).Range.Text = ""
$sec.Headers(3).Range.Text = ""
$sec.Footers(1).Range.Text = ""
$sec.Footers(2).Range.Text = ""
$sec.Footers(3).Range.Text = ""

Write-Output $d.Content.Text
